$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Datos actualizados" timestamp string (row 1)
$ws.Range("A1").Value = "Datos actualizados a 19 de Agosto de 2020 a las 16:01"

# Row 4
$ws.Cells.Item(4, 2).Value = 5657561
$ws.Cells.Item(4, 3).Value = 1587
$ws.Cells.Item(4, 4).Value = 3012058
$ws.Cells.Item(4, 5).Value = 2470384
$ws.Cells.Item(4, 6).Value = 0
$ws.Cells.Item(4, 7).Value = 45
$ws.Cells.Item(4, 8).Value = 175119

# Row 6
$ws.Cells.Item(6, 2).Value = 2786999
$ws.Cells.Item(6, 3).Value = 20373
$ws.Cells.Item(6, 4).Value = 2052736
$ws.Cells.Item(6, 5).Value = 681099
$ws.Cells.Item(6, 6).Value = 0
$ws.Cells.Item(6, 7).Value = 150
$ws.Cells.Item(6, 8).Value = 53164

# Row 16
$ws.Cells.Item(16, 2).Value = 305966
$ws.Cells.Item(16, 3).Value = 0
$ws.Cells.Item(16, 4).Value = 228725
$ws.Cells.Item(16, 5).Value = 71127
$ws.Cells.Item(16, 6).Value = 0
$ws.Cells.Item(16, 7).Value = 66
$ws.Cells.Item(16, 8).Value = 6114

# Row 17
$ws.Cells.Item(17, 2).Value = 302686
$ws.Cells.Item(17, 3).Value = 1363
$ws.Cells.Item(17, 4).Value = 274091
$ws.Cells.Item(17, 5).Value = 25089
$ws.Cells.Item(17, 6).Value = 0
$ws.Cells.Item(17, 7).Value = 36
$ws.Cells.Item(17, 8).Value = 3506

# Row 22
$ws.Cells.Item(22, 2).Value = 228138
$ws.Cells.Item(22, 3).Value = 33
$ws.Cells.Item(22, 4).Value = 203900
$ws.Cells.Item(22, 5).Value = 14933
$ws.Cells.Item(22, 6).Value = 0
$ws.Cells.Item(22, 7).Value = 0
$ws.Cells.Item(22, 8).Value = 9305

# Row 24
$ws.Cells.Item(24, 2).Value = 188802
$ws.Cells.Item(24, 3).Value = 4093
$ws.Cells.Item(24, 4).Value = 134369
$ws.Cells.Item(24, 5).Value = 48312
$ws.Cells.Item(24, 6).Value = 0
$ws.Cells.Item(24, 7).Value = 85
$ws.Cells.Item(24, 8).Value = 6121

# Row 28
$ws.Cells.Item(28, 2).Value = 115956
$ws.Cells.Item(28, 3).Value = 295
$ws.Cells.Item(28, 4).Value = 112658
$ws.Cells.Item(28, 5).Value = 3105
$ws.Cells.Item(28, 6).Value = 0
$ws.Cells.Item(28, 7).Value = 0
$ws.Cells.Item(28, 8).Value = 193

# Row 63
$ws.Cells.Item(63, 2).Value = 34620
$ws.Cells.Item(63, 3).Value = 146
$ws.Cells.Item(63, 4).Value = 32363
$ws.Cells.Item(63, 5).Value = 1747
$ws.Cells.Item(63, 6).Value = 0
$ws.Cells.Item(63, 7).Value = 1
$ws.Cells.Item(63, 8).Value = 510

# Row 67
$ws.Cells.Item(67, 2).Value = 30048
$ws.Cells.Item(67, 3).Value = 158
$ws.Cells.Item(67, 4).Value = 27702
$ws.Cells.Item(67, 5).Value = 1662
$ws.Cells.Item(67, 6).Value = 0
$ws.Cells.Item(67, 7).Value = 3
$ws.Cells.Item(67, 8).Value = 684

# Row 76
$ws.Cells.Item(76, 2).Value = 17606
$ws.Cells.Item(76, 3).Value = 300
$ws.Cells.Item(76, 4).Value = 10312
$ws.Cells.Item(76, 5).Value = 7179
$ws.Cells.Item(76, 6).Value = 0
$ws.Cells.Item(76, 7).Value = 2
$ws.Cells.Item(76, 8).Value = 115

# Row 78
$ws.Cells.Item(78, 2).Value = 16691
$ws.Cells.Item(78, 3).Value = 340
$ws.Cells.Item(78, 4).Value = 10711
$ws.Cells.Item(78, 5).Value = 5473
$ws.Cells.Item(78, 6).Value = 0
$ws.Cells.Item(78, 7).Value = 12
$ws.Cells.Item(78, 8).Value = 507

# Row 83
$ws.Cells.Item(83, 2).Value = 13076
$ws.Cells.Item(83, 3).Value = 106
$ws.Cells.Item(83, 4).Value = 9625
$ws.Cells.Item(83, 5).Value = 2900
$ws.Cells.Item(83, 6).Value = 0
$ws.Cells.Item(83, 7).Value = 2
$ws.Cells.Item(83, 8).Value = 551

# Row 87 -> country swapped to Zambia
$ws.Cells.Item(87, 1).Value = "Zambia"
$ws.Cells.Item(87, 2).Value = 10218
$ws.Cells.Item(87, 3).Value = 237
$ws.Cells.Item(87, 4).Value = 9126
$ws.Cells.Item(87, 5).Value = 823
$ws.Cells.Item(87, 6).Value = 0
$ws.Cells.Item(87, 7).Value = 5
$ws.Cells.Item(87, 8).Value = 269

# Row 88 -> country swapped to Noruega
$ws.Cells.Item(88, 1).Value = "Noruega"
$ws.Cells.Item(88, 2).Value = 10111
$ws.Cells.Item(88, 3).Value = 0
$ws.Cells.Item(88, 4).Value = 8857
$ws.Cells.Item(88, 5).Value = 992
$ws.Cells.Item(88, 6).Value = 0
$ws.Cells.Item(88, 7).Value = 0
$ws.Cells.Item(88, 8).Value = 262

# Row 96
$ws.Cells.Item(96, 2).Value = 8166
$ws.Cells.Item(96, 3).Value = 35
$ws.Cells.Item(96, 4).Value = 6971
$ws.Cells.Item(96, 5).Value = 1130
$ws.Cells.Item(96, 6).Value = 0
$ws.Cells.Item(96, 7).Value = 0
$ws.Cells.Item(96, 8).Value = 65

# Row 98
$ws.Cells.Item(98, 2).Value = 7805
$ws.Cells.Item(98, 3).Value = 29
$ws.Cells.Item(98, 4).Value = 7100
$ws.Cells.Item(98, 5).Value = 371
$ws.Cells.Item(98, 6).Value = 0
$ws.Cells.Item(98, 7).Value = 0
$ws.Cells.Item(98, 8).Value = 334

# Row 119
$ws.Cells.Item(119, 2).Value = 3265
$ws.Cells.Item(119, 3).Value = 8
$ws.Cells.Item(119, 4).Value = 2396
$ws.Cells.Item(119, 5).Value = 776
$ws.Cells.Item(119, 6).Value = 0
$ws.Cells.Item(119, 7).Value = 0
$ws.Cells.Item(119, 8).Value = 93

# Row 141 -> country swapped to Uganda
$ws.Cells.Item(141, 1).Value = "Uganda"
$ws.Cells.Item(141, 2).Value = 1656
$ws.Cells.Item(141, 3).Value = 53
$ws.Cells.Item(141, 4).Value = 1188
$ws.Cells.Item(141, 5).Value = 452
$ws.Cells.Item(141, 6).Value = 0
$ws.Cells.Item(141, 7).Value = 1
$ws.Cells.Item(141, 8).Value = 16

# Row 142 -> country swapped to Nueva Zelanda
$ws.Cells.Item(142, 1).Value = "Nueva Zelanda"
$ws.Cells.Item(142, 2).Value = 1649
$ws.Cells.Item(142, 3).Value = 6
$ws.Cells.Item(142, 4).Value = 1531
$ws.Cells.Item(142, 5).Value = 96
$ws.Cells.Item(142, 6).Value = 0
$ws.Cells.Item(142, 7).Value = 0
$ws.Cells.Item(142, 8).Value = 22

# Row 184
$ws.Cells.Item(184, 2).Value = 223
$ws.Cells.Item(184, 3).Value = 1
$ws.Cells.Item(184, 4).Value = 197
$ws.Cells.Item(184, 5).Value = 26
$ws.Cells.Item(184, 6).Value = 0
$ws.Cells.Item(184, 7).Value = 0
$ws.Cells.Item(184, 8).Value = 0

# Row 213 -> country swapped to Montserrat
$ws.Cells.Item(213, 1).Value = "Montserrat"
$ws.Cells.Item(213, 2).Value = 13
$ws.Cells.Item(213, 3).Value = 0
$ws.Cells.Item(213, 4).Value = 12
$ws.Cells.Item(213, 5).Value = 0
$ws.Cells.Item(213, 6).Value = 0
$ws.Cells.Item(213, 7).Value = 0
$ws.Cells.Item(213, 8).Value = 1

# Row 214 -> country swapped to Islas Malvinas
$ws.Cells.Item(214, 1).Value = "Islas Malvinas"
$ws.Cells.Item(214, 2).Value = 13
$ws.Cells.Item(214, 3).Value = 0
$ws.Cells.Item(214, 4).Value = 13
$ws.Cells.Item(214, 5).Value = 0
$ws.Cells.Item(214, 6).Value = 0
$ws.Cells.Item(214, 7).Value = 0
$ws.Cells.Item(214, 8).Value = 0
